$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 244, shifting existing rows 244:290 down to 245:291
$ws.Rows.Item(244).Insert()

# Populate the newly inserted row 244 with the new record's data
$ws.Cells.Item(244, 1).Value = 9
$ws.Cells.Item(244, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(244, 3).Value = "Metropolitana"
$ws.Cells.Item(244, 4).Value = 44694
$ws.Cells.Item(244, 5).Value = 13
$ws.Cells.Item(244, 6).Value = 100112021
$ws.Cells.Item(244, 7).Value = "Ají"
$ws.Cells.Item(244, 8).Value = "Americana (o)"
$ws.Cells.Item(244, 9).Value = "Primera"
$ws.Cells.Item(244, 10).Value = 16
$ws.Cells.Item(244, 11).Value = 25000
$ws.Cells.Item(244, 12).Value = 26000
$ws.Cells.Item(244, 13).Value = 25500
$ws.Cells.Item(244, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(244, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(244, 16).Value = 1020
$ws.Cells.Item(244, 17).Value = 25
$ws.Cells.Item(244, 18).Value = "Hortaliza"

Write-Host "done"
